# betavals.xlsx update:
#  - Insert 6 new ticker rows (GBP.JPY, GLD, HG, SD, SPCE, VALE) at their
#    alphabetically-sorted positions in the existing "Ticker" list.
#  - All pre-existing rows (and the two formulas that reference other rows
#    by address) shift down automatically with Rows.Insert(); Excel keeps
#    formulas pointing at the *same ticker* they pointed at before.
#  - Fill the values/comments for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new blank rows, bottom-most first, so the row numbers below
# (computed against the *original* layout) stay valid while we work.
$ws.Rows.Item(41).Insert()   # new row for VALE (before VNQ)
$ws.Rows.Item(38).Insert()   # new row for SPCE (before STNG)
$ws.Rows.Item(35).Insert()   # new row for SD   (before SHOP)
$ws.Rows.Item(22).Insert()   # new row for HG   (before HYG)
$ws.Rows.Item(19).Insert()   # new row for GLD  (before GLEN)
$ws.Rows.Item(16).Insert()   # new row for GBP.JPY (before GC)

# Fill in the new rows.
$ws.Range("A16").Value = "GBP.JPY"
$ws.Range("B16").Value = -0.0061
$ws.Range("D16").Value = "GBJP not right but value in JPY which is too big"

$ws.Range("A20").Value = "GLD"
$ws.Range("B20").Value = 0.02

$ws.Range("A24").Value = "HG"
$ws.Range("B24").Value = 0.72
$ws.Range("D24").Value = "spy vs jjc"

$ws.Range("A38").Value = "SD"
$ws.Range("B38").Value = 0.081

$ws.Range("A42").Value = "SPCE"
$ws.Range("B42").Value = 2

$ws.Range("A46").Value = "VALE"
$ws.Range("B46").Value = 0.43

# Match the author's last-saved selection.
$ws.Range("D16").Select()
